# Update song-list worksheet to reflect the 2026/2/23 listening session.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Cells.Item(23, 5).Value = "2025/03/04，2025/05/01，2025/05/08，2025/05/26，2025/06/05，2025/06/23，2025/06/26，2025/07/03，2025/07/20，2025/07/26，2025/08/03，2025/08/12，2025/08/30，2025/09/21，2025/09/29，2025/11/9，2025/11/12，2025/11/27，2025/12/4，2025/12/28，2026/1/22，2026/2/2，2026/2/23"
$ws.Cells.Item(23, 8).Value = "'23"

# Row 24
$ws.Cells.Item(24, 5).Value = "2025/03/04，2025/04/14，2025/05/09，2025/06/13，2025/06/19，2025/06/27，2025/07/10，2025/09/10，2025/09/19，2025/11/14，2025/12/31，2026/2/23"
$ws.Cells.Item(24, 8).Value = "'12"
$ws.Cells.Item(24, 9).Value = "BV1r4oLYfEjE，BV1Asf6BvEo3"

# Row 27
$ws.Cells.Item(27, 5).Value = "2025/03/04，2025/03/14，2025/04/15，2025/05/03，2025/05/19，2025/07/18，2025/07/28，2025/08/28，2025/09/22，2025/10/19，2025/10/26，2025/11/9，2025/12/30，2026/2/20，2026/2/23"
$ws.Cells.Item(27, 8).Value = "'15"

# Row 60
$ws.Cells.Item(60, 5).Value = "2025/03/07，2025/03/17，2025/03/28，2025/03/29，2025/04/13，2025/04/15，2025/04/18，2025/04/28，2025/05/22，2025/06/02，2025/06/03，2025/06/10，2025/06/12，2025/06/20，2025/06/20，2025/06/27，2025/07/05，2025/07/09，2025/07/13，2025/08/04，2025/08/08，2025/08/12，2025/08/24，2025/09/06，2025/10/07，2025/10/24，2025/11/24，2025/12/4，2025/12/24，2026/1/15，2026/2/23"
$ws.Cells.Item(60, 8).Value = "'31"

# Row 61
$ws.Cells.Item(61, 5).Value = "2025/03/07，2025/04/04，2025/04/12，2025/04/18，2025/04/24，2025/06/17，2025/06/19，2025/07/11，2025/07/18，2025/08/28，2025/09/05，2025/09/08，2025/09/27，2025/10/03，2025/12/11，2025/12/30，2026/1/21，2026/1/25，2026/2/10，2026/2/23"
$ws.Cells.Item(61, 8).Value = "'20"

# Row 64
$ws.Cells.Item(64, 5).Value = "2025/03/07，2025/04/11，2025/04/29，2025/08/01，2025/10/26，2026/2/19，2026/2/23"
$ws.Cells.Item(64, 8).Value = "'7"

# Row 71
$ws.Cells.Item(71, 5).Value = "2025/03/09，2025/03/14，2025/04/07，2025/04/15，2025/04/27，2025/05/12，2025/05/27，2025/06/05，2025/06/08，2025/07/20，2025/08/01，2025/09/21，2025/09/28，2025/10/07，2025/10/18，2025/10/24，2025/11/10，2025/12/1，2026/1/16，2026/2/2，2026/2/23"
$ws.Cells.Item(71, 8).Value = "'21"

# Row 130
$ws.Cells.Item(130, 5).Value = "2025/03/17，2025/03/28，2025/04/15，2025/04/19，2025/04/25，2025/05/01，2025/05/18，2025/05/25，2025/05/29，2025/06/05，2025/06/12，2025/06/22，2025/07/03，2025/07/11，2025/07/18，2025/07/22，2025/07/31，2025/08/01，2025/08/02，2025/08/03，2025/08/12，2025/09/03，2025/09/11，2025/09/18，2025/09/26，2025/10/03，2025/10/09，2025/10/20，2025/11/12，2025/12/7，2025/12/23，2026/1/5，2026/1/21，2026/1/29，2026/2/9，2026/2/20，2026/2/23"
$ws.Cells.Item(130, 8).Value = "'37"
$ws.Cells.Item(130, 9).Value = "BV1A1fzBwEjQ，BV19KmPBDERd，BV1B1fGBjEtL"

# Row 153
$ws.Cells.Item(153, 5).Value = "2025/03/20，2025/05/19，2025/06/12，2025/06/19，2025/06/23，2025/07/11，2025/08/03，2025/08/12，2025/08/21，2025/09/15，2025/10/07，2025/12/2，2026/1/16，2026/2/2，2026/2/23"
$ws.Cells.Item(153, 8).Value = "'15"

# Row 157
$ws.Cells.Item(157, 5).Value = "2025/03/20，2025/04/04，2025/05/03，2025/05/22，2025/06/05，2025/07/25，2025/08/26，2025/08/16，2025/09/21，2025/10/03，2025/11/12，2026/1/23，2026/1/26，2026/2/7，2026/2/10，2026/2/23"
$ws.Cells.Item(157, 8).Value = "'16"

# Row 184
$ws.Cells.Item(184, 5).Value = "2025/03/22，2025/04/05，2025/04/17，2025/04/21，2025/04/22，2025/04/25，2025/05/08，2025/05/15，2025/05/25，2025/05/29，2025/06/03，2025/06/13，2025/06/23，2025/07/06，2025/07/15，2025/08/03，2025/08/16，2025/09/03，2025/09/10，2025/09/30，2025/10/16，2025/10/22，2025/10/30，2025/11/18，2025/12/2，2025/12/5，2025/12/30，2026/1/23，2026/1/29，2026/2/23"
$ws.Cells.Item(184, 8).Value = "'30"
$ws.Cells.Item(184, 9).Value = "BV1Xr2xBVEqN，BV1Csf6BeE7B"

# Row 208
$ws.Cells.Item(208, 5).Value = "2025/03/27，2025/04/03，2025/04/15，2025/04/24，2025/05/13，2025/06/02，2025/06/27，2025/07/14，2025/08/26，2025/09/02，2025/10/03，2025/10/17，2025/11/21，2026/2/14，2026/2/23"
$ws.Cells.Item(208, 8).Value = "'15"
$ws.Cells.Item(208, 9).Value = "BV1Z5UEBvEPF，BV12Ef6BBEgR"

# Row 301
$ws.Cells.Item(301, 5).Value = "2025/04/06，2025/04/19，2025/05/03，2025/08/01，2025/10/31，2025/12/15，2026/2/22，2026/2/23"
$ws.Cells.Item(301, 8).Value = "'8"

# Row 312
$ws.Cells.Item(312, 5).Value = "2025/04/08，2025/05/12，2025/08/22，2025/10/01，2025/10/17，2026/2/23"
$ws.Cells.Item(312, 8).Value = "'6"
$ws.Cells.Item(312, 9).Value = "BV1cnWpzrEtD，BV1FSfzBbEpH"

# Row 340
$ws.Cells.Item(340, 5).Value = "2025/04/13，2025/07/13，2025/07/18，2025/08/02，2025/09/04，2025/09/12，2025/10/17，2025/11/28，2026/1/22，2026/2/23"
$ws.Cells.Item(340, 8).Value = "'10"
$ws.Cells.Item(340, 9).Value = "BV1GBaQzcE39"

# Row 355
$ws.Cells.Item(355, 5).Value = "2025/04/17，2025/05/01，2025/05/26，2025/07/03，2025/07/27，2025/09/18，2025/10/17，2025/11/21，2026/1/12，2026/1/24，2026/2/19，2026/2/23"
$ws.Cells.Item(355, 8).Value = "'12"
$ws.Cells.Item(355, 9).Value = "BV13E33zFES7，BV1Knf6BMELa"

# Row 442
$ws.Cells.Item(442, 9).Value = "BV1J67NzGEDg，BV173fNB4EZX"

# Row 498
$ws.Cells.Item(498, 5).Value = "2025/06/01，2026/1/23，2026/2/2，2026/2/23"
$ws.Cells.Item(498, 8).Value = "'4"
$ws.Cells.Item(498, 9).Value = "BV1w2fzBhEHL，BV1VuzuB1EQf"

# Row 587
$ws.Cells.Item(587, 5).Value = "2025/07/14，2025/07/17，2025/07/18，2025/07/24，2025/07/28，2025/08/02，2025/08/03，2025/08/21，2025/09/05，2025/09/18，2025/09/30，2025/10/19，2025/11/3，2025/11/22，2026/1/5，2026/2/18，2026/2/23"
$ws.Cells.Item(587, 8).Value = "'17"

# Row 631
$ws.Cells.Item(631, 5).Value = "2025/09/10，2025/09/22，2025/09/24，2025/10/23，2025/12/6，2026/1/24，2026/2/5，2026/2/23"
$ws.Cells.Item(631, 8).Value = "'8"

# Row 650
$ws.Cells.Item(650, 5).Value = "2025/09/17，2025/10/11，2025/10/13，2025/10/20，2025/10/25，2025/11/21，2025/12/17，2026/2/23"
$ws.Cells.Item(650, 8).Value = "'8"
$ws.Cells.Item(650, 9).Value = "BV1aZWyzREbb，BV1mLf6B6EYy"

# Row 717
$ws.Cells.Item(717, 5).Value = "2025/10/24，2025/10/25，2025/10/28，2025/11/12，2025/11/26，2026/1/2，2026/1/5，2026/1/9，2026/2/23"
$ws.Cells.Item(717, 8).Value = "'9"
$ws.Cells.Item(717, 9).Value = "BV1TBihBsEo5，BV1ksf6BeE2k"

# Row 719
$ws.Cells.Item(719, 5).Value = "2025/10/25，2025/12/2，2026/2/23"
$ws.Cells.Item(719, 8).Value = "'3"
$ws.Cells.Item(719, 9).Value = "BV1oV2iBGEn4"

# Row 769
$ws.Cells.Item(769, 5).Value = "2025/12/19，2025/12/23，2025/12/26，2026/1/2，2026/1/5，2026/1/11，2026/1/25，2026/2/9，2026/2/23"
$ws.Cells.Item(769, 8).Value = "'9"
$ws.Cells.Item(769, 9).Value = "BV1oKiwBXE1L，BV1Csf6BeEYe"

# Row 773
$ws.Cells.Item(773, 5).Value = "2025/12/26，2026/1/10，2026/1/22，2026/1/23，2026/1/31，2026/2/10，2026/2/23"
$ws.Cells.Item(773, 8).Value = "'7"

# Row 788
$ws.Cells.Item(788, 5).Value = "2026/1/25，2026/1/30，2026/2/2，2026/2/8，2026/2/19，2026/2/23"
$ws.Cells.Item(788, 8).Value = "'6"

# Row 796
$ws.Cells.Item(796, 5).Value = "2026/2/8，2026/2/9，2026/2/12，2026/2/14，2026/2/19，2026/2/23"
$ws.Cells.Item(796, 8).Value = "'6"
$ws.Cells.Item(796, 9).Value = "BV1W4cuzSE8Y，BV12Ef6BBE5C"

# Row 806
$ws.Cells.Item(806, 5).Value = "2026/2/22，2026/2/23"
$ws.Cells.Item(806, 8).Value = "'2"
$ws.Cells.Item(806, 9).Value = "BV1K7f6BgE9r"

# New rows appended at the end (2026/2/23 session)
# Row 807
$ws.Cells.Item(807, 2).Value = "安静"
$ws.Cells.Item(807, 4).Value = "周杰伦"
$ws.Cells.Item(807, 5).Value = "'2026/2/23"
$ws.Cells.Item(807, 7).Value = "华语"
$ws.Cells.Item(807, 8).Value = "'1"
$ws.Cells.Item(807, 9).Value = "BV11nf6BTE4b"

# Row 808
$ws.Cells.Item(808, 2).Value = "可惜没如果"
$ws.Cells.Item(808, 4).Value = "林俊杰"
$ws.Cells.Item(808, 5).Value = "'2026/2/23"
$ws.Cells.Item(808, 7).Value = "华语"
$ws.Cells.Item(808, 8).Value = "'1"
$ws.Cells.Item(808, 9).Value = "BV1w1fzBwEhj"
